$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ETH")
$ws.Range("J3").Value = 2772.397703171068
$ws.Range("B12").Value = 0.00725439
$ws.Range("B36").Value = 0.02517439
$ws.Range("D36").Value = 44.3
$ws.Range("B40").Value = 0.05707998
$ws.Range("D40").Value = 107.45
$ws = $wb.Worksheets.Item("APE")
$ws.Range("J3").Value = 1.603596163217528
$ws.Range("B5").Value = 16.23066542
$ws.Range("D5").Value = 44.3
$ws.Range("B6").Value = 0.59997618
$ws = $wb.Worksheets.Item("ATOM")
$ws.Range("J3").Value = 10.13462526122557
$ws.Range("B7").Value = 0.02909975
$ws = $wb.Worksheets.Item("AVAX")
$ws.Range("J3").Value = 39.51908498948486
$ws.Range("B5").Value = 2.65463808
$ws.Range("D5").Value = 44.3
$ws.Range("B6").Value = 0.0167139
$ws = $wb.Worksheets.Item("AMP")
$ws.Range("J3").Value = 0.003712244187653699
$ws = $wb.Worksheets.Item("BNB")
$ws.Range("J3").Value = 358.5720806687447
$ws.Range("B10").Value = 0.00272536
$ws.Range("B12").Value = 0.1544037
$ws.Range("D12").Value = 44.3
$ws = $wb.Worksheets.Item("DOGE")
$ws.Range("J3").Value = 0.0841828193576163
$ws.Range("B6").Value = 0.29120616
$ws = $wb.Worksheets.Item("DOT")
$ws.Range("J3").Value = 7.590810143631857
$ws.Range("B5").Value = 7.73572629
$ws.Range("D5").Value = 44.3
$ws.Range("B6").Value = 0.07958381
$ws = $wb.Worksheets.Item("EGLD")
$ws.Range("J3").Value = 57.41993011092561
$ws.Range("B6").Value = 0.00299655
$ws = $wb.Worksheets.Item("GRT")
$ws.Range("J3").Value = 0.1878736433047423
$ws = $wb.Worksheets.Item("ICP")
$ws.Range("J3").Value = 13.10717360326807
$ws.Range("B6").Value = 0.00235864
$ws = $wb.Worksheets.Item("BTC")
$ws.Range("J3").Value = 51598.29631906211
$ws.Range("B6").Value = 0.00035572
$ws.Range("B24").Value = 0.00165298
$ws.Range("D24").Value = 44.3
$ws.Range("B34").Value = 0.00208117
$ws.Range("D34").Value = 61.55
$ws = $wb.Worksheets.Item("KAVA")
$ws.Range("J3").Value = 0.7392583809834428
$ws = $wb.Worksheets.Item("LDO")
$ws.Range("J3").Value = 3.127978494216071
$ws.Range("B6").Value = 0.02019336
$ws = $wb.Worksheets.Item("LINK")
$ws.Range("J3").Value = 19.63423207165748
$ws.Range("B6").Value = 0.00248108
$ws = $wb.Worksheets.Item("LTC")
$ws.Range("J3").Value = 69.81431374808243
$ws.Range("B6").Value = 0.00133658
$ws = $wb.Worksheets.Item("LUNA")
$ws.Range("J3").Value = 0.7055846520645814
$ws.Range("B6").Value = 0.05840602
$ws = $wb.Worksheets.Item("LUNC")
$ws.Range("J3").Value = 0.0001295553757483842
$ws.Range("B18").Value = 5023.23493174
$ws = $wb.Worksheets.Item("MATIC")
$ws.Range("J3").Value = 0.9231072342724557
$ws.Range("B6").Value = 0.3290304
$ws.Range("B7").Value = 48.74884497
$ws.Range("D7").Value = 44.3
$ws = $wb.Worksheets.Item("MEME")
$ws.Range("J3").Value = 0.02589042453729953
$ws.Range("B6").Value = 0.06773084
$ws = $wb.Worksheets.Item("MINA")
$ws.Range("J3").Value = 1.346369914660285
$ws.Range("B6").Value = 0.34996814
$ws = $wb.Worksheets.Item("NEAR")
$ws.Range("J3").Value = 3.240298516082145
$ws.Range("B6").Value = 23.96020479
$ws.Range("D6").Value = 44.3
$ws.Range("B7").Value = 0.10294619
$ws = $wb.Worksheets.Item("SEI")
$ws.Range("J3").Value = 0.9262281038994522
$ws.Range("B6").Value = 0.07635641999999999
$ws = $wb.Worksheets.Item("SHIB")
$ws.Range("J3").Value = 0.00000963982721493776
$ws.Range("B6").Value = 275.69
$ws = $wb.Worksheets.Item("SHPING")
$ws.Range("J3").Value = 0.004829254508509664
$ws = $wb.Worksheets.Item("SOL")
$ws.Range("J3").Value = 108.1519584255934
$ws.Range("B17").Value = 0.06471251
$ws.Range("B18").Value = 1.9171108
$ws.Range("D18").Value = 44.3
$ws = $wb.Worksheets.Item("TRX")
$ws.Range("J3").Value = 0.1380240610495743
$ws.Range("B6").Value = 0.26513448
$ws = $wb.Worksheets.Item("UNI")
$ws.Range("J3").Value = 7.654681341834773
$ws.Range("B6").Value = 0.0027469
$ws = $wb.Worksheets.Item("XRP")
$ws.Range("J3").Value = 0.5526457767990984
$ws.Range("B6").Value = 0.86487264
$ws = $wb.Worksheets.Item("TIA")
$ws.Range("J3").Value = 18.66235392418715
$ws.Range("B6").Value = 0.00404968
$ws = $wb.Worksheets.Item("DYDX")
$ws.Range("J3").Value = 3.039948510596013
$ws.Range("B6").Value = 0.00080303
$ws = $wb.Worksheets.Item("POLIS")
$ws.Range("J3").Value = 0.4319540847869544
$ws = $wb.Worksheets.Item("ATLAS")
$ws.Range("J3").Value = 0.004821947522984643
$ws = $wb.Worksheets.Item("ACE")
$ws.Range("J3").Value = 9.808249599470496
$ws.Range("B6").Value = 0.00002495
$ws = $wb.Worksheets.Item("ADA")
$ws.Range("J3").Value = 0.5869308421824619
$ws.Range("B6").Value = 0.7811301899999999
$ws.Range("B7").Value = 122.82958283
$ws.Range("D7").Value = 44.3
$ws = $wb.Worksheets.Item("ALGO")
$ws.Range("J3").Value = 0.1900243911060043
$ws.Range("B6").Value = 0.58418022
